Write-Host "ok"
